$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking Price strings from Excel auto-number-conversion
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.293.85'
$ws.Range('D3').Value = '1.668.62'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('D5').Value = '218.68'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').Value = '0.5245'
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('D9').Value = '0.06334'
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('E10').Value = '  -3.17%  '
$ws.Range('D11').Value = '0.07766'
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').Value = '1.677.78'
$ws.Range('D13').Value = '4.455'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').Value = '1.894.91'
$ws.Range('E14').Value = '  -0.67%  '
$ws.Range('D15').Value = '0.5500'
$ws.Range('E15').Value = '  -1.18%  '
$ws.Range('D16').Value = '0.0₅8288'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = '65.07'
$ws.Range('E17').Value = '  -0.84%  '
$ws.Range('D18').Value = '26.322.23'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').Value = '4.686'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').Value = '195.06'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').Value = '10.16'
$ws.Range('E22').Value = '  -1.27%  '
$ws.Range('D23').Value = '6.084'
$ws.Range('E23').Value = '  -4.04%  '
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').Value = '139.92'
$ws.Range('E25').Value = '  -1.61%  '
$ws.Range('D26').Value = '0.1240'
$ws.Range('E26').Value = '  -3.72%  '
$ws.Range('D27').Value = '7.214'
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('D28').Value = '16.20'
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('D29').Value = '1.416'
$ws.Range('E29').Value = '  -1.51%  '
$ws.Range('D30').Value = '0.06182'
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('D31').Value = '1.283'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D32').Value = '3.593'
$ws.Range('E32').Value = '  -0.56%  '
$ws.Range('D33').Value = '3.302'
$ws.Range('E33').Value = '  -4.31%  '
$ws.Range('E34').Value = '  -2.64%  '
$ws.Range('D35').Value = '0.9740'
$ws.Range('E35').Value = '  -3.34%  '
$ws.Range('D36').Value = '2.428'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '2.793'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').Value = '0.5765'
$ws.Range('E38').Value = '  -5.83%  '
$ws.Range('D39').Value = '0.01609'
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('D40').Value = '6.022'
$ws.Range('E40').Value = '  -2.02%  '
$ws.Range('D41').Value = '0.8592'
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').Value = '1.024.72'
$ws.Range('E43').Value = '  -5.54%  '
$ws.Range('D44').Value = '100.37'
$ws.Range('E44').Value = '  +0.02%  '
$ws.Range('D45').Value = '1.810.86'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '57.78'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₈108'
$ws.Range('E47').Value = '  +4.77%  '
$ws.Range('D48').Value = '1.007'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D49').Value = '8.061'
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('D50').Value = '1.489'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').Value = '0.05183'
$ws.Range('E51').Value = '  -0.49%  '

# Restore default (no explicit) style on the Price column so no stray formatting remains
$ws.Range("D2:D51").Style = "Normal"
